# Apply updated crypto price / 1h volume data to the sheet.
# Prices in column D are numeric-looking text (e.g. "55.935.57"), so a
# leading apostrophe is used to force Excel to store them as text, exactly
# like the source data (matches original inlineStr string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'55.935.57"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("D3").Value = "'2.360.14"
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'500.54"
$ws.Range("E5").Value = "  -2.10%  "

$ws.Range("D6").Value = "'129.11"
$ws.Range("E6").Value = "  -3.38%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("E8").Value = "  -3.00%  "

$ws.Range("D9").Value = "'2.359.76"
$ws.Range("E9").Value = "  -4.11%  "

$ws.Range("D10").Value = "'0.0979"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'4.78"
$ws.Range("E12").Value = "  +2.96%  "

$ws.Range("D13").Value = "'0.323"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").Value = "'2.777.01"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("D15").Value = "'55.912.05"
$ws.Range("E15").Value = "  -3.34%  "

$ws.Range("D16").Value = "'21.36"
$ws.Range("E16").Value = "  -2.94%  "

$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").Value = "'2.416.79"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").Value = "'9.99"
$ws.Range("E19").Value = "  -3.70%  "

$ws.Range("E20").Value = "  -3.76%  "

$ws.Range("D21").Value = "'306.53"
$ws.Range("E21").Value = "  -2.81%  "

$ws.Range("D22").Value = "'6.26"
$ws.Range("E22").Value = "  -3.20%  "

$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "'65.79"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").Value = "'0.368"
$ws.Range("E26").Value = "  -3.96%  "

$ws.Range("E27").Value = "  -6.48%  "

$ws.Range("D28").Value = "'7.21"
$ws.Range("E28").Value = "  -5.48%  "

$ws.Range("D29").Value = "'172.05"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("E30").Value = "  -4.22%  "

$ws.Range("D31").Value = "'1.64"
$ws.Range("E31").Value = "  -3.33%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'5.75"
$ws.Range("E33").Value = "  -6.69%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("D35").Value = "'1.08"
$ws.Range("E35").Value = "  -6.34%  "

$ws.Range("D36").Value = "'17.56"
$ws.Range("E36").Value = "  -3.05%  "

$ws.Range("D37").Value = "'1.17"
$ws.Range("E37").Value = "  -6.26%  "

$ws.Range("D38").Value = "'3.71"
$ws.Range("E38").Value = "  -5.13%  "

$ws.Range("D39").Value = "'36.13"
$ws.Range("E39").Value = "  -1.71%  "

$ws.Range("D40").Value = "'0.791"
$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("D41").Value = "'1.38"
$ws.Range("E41").Value = "  -6.77%  "

$ws.Range("D42").Value = "'3.35"
$ws.Range("E42").Value = "  -2.06%  "

$ws.Range("D43").Value = "'128.51"
$ws.Range("E43").Value = "  -6.10%  "

$ws.Range("D44").Value = "'4.68"
$ws.Range("E44").Value = "  -6.96%  "

$ws.Range("D45").Value = "'0.561"
$ws.Range("E45").Value = "  -3.16%  "

$ws.Range("D46").Value = "'0.0899"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("D47").Value = "'237.42"
$ws.Range("E47").Value = "  -7.16%  "

$ws.Range("D48").Value = "'0.0480"
$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D49").Value = "'0.0206"
$ws.Range("E49").Value = "  -4.44%  "

$ws.Range("D50").Value = "'16.95"
$ws.Range("E50").Value = "  -2.79%  "

$ws.Range("D51").Value = "'0.948"
$ws.Range("E51").Value = "  -1.28%  "
